$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.333.91'
$ws.Cells.Item(2, 5).Value = '  +12.56%  '

$ws.Cells.Item(3, 4).Value = '1.825.19'
$ws.Cells.Item(3, 5).Value = '  +9.25%  '

$ws.Cells.Item(4, 4).Value = '''0.998'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.26%  '

$ws.Cells.Item(5, 4).Value = '''229.69'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +4.70%  '

$ws.Cells.Item(6, 5).Value = '  +4.20%  '

$ws.Cells.Item(7, 4).Value = '''0.998'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.34%  '

$ws.Cells.Item(8, 5).Value = '  +6.90%  '

$ws.Cells.Item(9, 4).Value = '''47.37'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +6.85%  '

$ws.Cells.Item(10, 4).Value = '''0.284'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +7.80%  '

$ws.Cells.Item(11, 5).Value = '  +5.97%  '

$ws.Cells.Item(12, 5).Value = '  +2.88%  '

$ws.Cells.Item(13, 4).Value = '2.088.75'
$ws.Cells.Item(13, 5).Value = '  +9.36%  '

$ws.Cells.Item(14, 4).Value = '1.827.48'
$ws.Cells.Item(14, 5).Value = '  +9.31%  '

$ws.Cells.Item(15, 5).Value = '  +5.85%  '

$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value = '34.281.29'
$ws.Cells.Item(16, 5).Value = '  +12.38%  '

$ws.Cells.Item(17, 2).Value = 'Chainlink'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(17, 4).Value = '''10.39'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +3.33%  '

$ws.Cells.Item(18, 4).Value = '''4.28'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +8.42%  '

$ws.Cells.Item(19, 4).Value = '''69.90'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +5.82%  '

$ws.Cells.Item(20, 4).Value = '''259.29'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +7.10%  '

$ws.Cells.Item(22, 4).Value = '''0.998'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.35%  '

$ws.Cells.Item(23, 5).Value = '  +6.85%  '

$ws.Cells.Item(24, 5).Value = '  +2.78%  '

$ws.Cells.Item(25, 5).Value = '  +3.59%  '

$ws.Cells.Item(26, 4).Value = '''158.80'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.48%  '

$ws.Cells.Item(27, 4).Value = '''16.70'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +5.81%  '

$ws.Cells.Item(28, 5).Value = '  +7.52%  '

$ws.Cells.Item(29, 5).Value = '  +2.58%  '

$ws.Cells.Item(30, 5).Value = '  -0.35%  '

$ws.Cells.Item(31, 5).Value = '  +12.64%  '

$ws.Cells.Item(32, 4).Value = '''0.0516'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +4.80%  '

$ws.Cells.Item(33, 5).Value = '  +6.20%  '

$ws.Cells.Item(34, 4).Value = '''3.55'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +8.63%  '

$ws.Cells.Item(35, 4).Value = '1.552.75'
$ws.Cells.Item(35, 5).Value = '  +3.66%  '

$ws.Cells.Item(36, 5).Value = '  +4.44%  '

$ws.Cells.Item(37, 5).Value = '  +7.01%  '

$ws.Cells.Item(38, 4).Value = '''85.69'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +2.53%  '

$ws.Cells.Item(39, 4).Value = '''0.633'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +7.22%  '

$ws.Cells.Item(40, 5).Value = '  +5.66%  '

$ws.Cells.Item(41, 5).Value = '  +5.86%  '

$ws.Cells.Item(42, 5).Value = '  +10.39%  '

$ws.Cells.Item(43, 4).Value = '''2.34'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +1.23%  '

$ws.Cells.Item(44, 5).Value = '  +10.50%  '

$ws.Cells.Item(45, 5).Value = '  +5.88%  '

$ws.Cells.Item(46, 2).Value = 'MinaProtocolToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Cells.Item(46, 4).Value = '''1.11'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +171.27%  '

$ws.Cells.Item(47, 2).Value = 'WEMIXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(47, 4).Value = '''1.07'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +4.92%  '

$ws.Cells.Item(48, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(48, 4).Value = '1.988.56'
$ws.Cells.Item(48, 5).Value = '  +10.24%  '

$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).Value = '''12.09'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +26.04%  '

$ws.Cells.Item(50, 2).Value = 'FraxShare'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(50, 4).Value = '''5.75'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +4.14%  '

$ws.Cells.Item(51, 2).Value = 'BitcoinSV'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Cells.Item(51, 4).Value = '''53.21'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +4.63%  '

